$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the two new sprint sheets *before* touching "2018.30.01" so the
#    Copy() snapshots its still-pristine (original) state.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2018.30.01")

# --- new sheet "2018.06.02" (sheetId 25) ------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$ws25 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws25.Name = "2018.06.02"

$ws25.Range("G3").Value = "(krank)"
$ws25.Range("C7").Value = 6
$ws25.Range("D7").Value = 6
$ws25.Range("E7").Value = 3
$ws25.Range("B8:E8").ClearContents()
$ws25.Range("A1:F16").Select()

# --- new sheet "Tabelle2" (sheetId 26) --------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$ws26 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws26.Name = "Tabelle2"

$ws26.Range("B3").Value = "Reviews & Statistics"
$ws26.Range("C3").Value = 3
$ws26.Range("D3").Value = 3
$ws26.Range("B4:F4").ClearContents()
$ws26.Range("B8:E8").ClearContents()
$ws26.Range("B14").Value = 3
$ws26.Range("B15").Value = 3
$ws26.Columns.Item(2).ColumnWidth = 15.592447916666666
$ws26.Columns.Item(3).ColumnWidth = 16.877604166666668
$ws26.Range("D5").Select()

# ---------------------------------------------------------------------------
# 2. Record the "worked" hours on the existing "2018.30.01" sheet.
# ---------------------------------------------------------------------------
$template.Range("E3").Value = 4
$template.Range("E4").Value = 2
$template.Range("E7").Value = 3
$template.Range("E8").Value = 3
$template.Range("A1:F16").Select()

# ---------------------------------------------------------------------------
# 3. Summary sheet: move the selection (it also loses tabSelected once a
#    later sheet becomes the active one).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("F24").Select()

# ---------------------------------------------------------------------------
# 4. Leave "Tabelle2" as the active sheet/tab.
# ---------------------------------------------------------------------------
$ws26.Activate()
$ws26.Range("D5").Select()
